$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing negative-word rows (3-37): words & stats shifted up/changed ---
# Row 4
$ws.Cells.Item(4, 2).Value = 0.8181818181818182
$ws.Cells.Item(4, 3).Value = 36
$ws.Cells.Item(4, 4).Value = 36
$ws.Cells.Item(4, 8).Value = 8
# Row 5
$ws.Cells.Item(5, 1).Value = "broke"
$ws.Cells.Item(5, 2).Value = 0.7475728155339806
$ws.Cells.Item(5, 3).Value = 154
$ws.Cells.Item(5, 4).Value = 154
$ws.Cells.Item(5, 8).Value = 52
# Row 7
$ws.Cells.Item(7, 1).Value = "however"
$ws.Cells.Item(7, 2).Value = 0.71875
$ws.Cells.Item(7, 3).Value = 46
$ws.Cells.Item(7, 4).Value = 46
$ws.Cells.Item(7, 8).Value = 18
# Row 8
$ws.Cells.Item(8, 1).Value = "disappointed"
$ws.Cells.Item(8, 2).Value = 0.6935483870967742
$ws.Cells.Item(8, 3).Value = 129
$ws.Cells.Item(8, 4).Value = 129
$ws.Cells.Item(8, 8).Value = 57
# Row 9
$ws.Cells.Item(9, 1).Value = "waste"
$ws.Cells.Item(9, 2).Value = 0.668918918918919
$ws.Cells.Item(9, 3).Value = 99
$ws.Cells.Item(9, 4).Value = 99
$ws.Cells.Item(9, 8).Value = 49
# Row 10
$ws.Cells.Item(10, 1).Value = "junk"
$ws.Cells.Item(10, 2).Value = 0.6363636363636364
$ws.Cells.Item(10, 3).Value = 35
$ws.Cells.Item(10, 4).Value = 35
$ws.Cells.Item(10, 8).Value = 20
# Row 11
$ws.Cells.Item(11, 1).Value = "guess"
$ws.Cells.Item(11, 2).Value = 0.6111111111111112
$ws.Cells.Item(11, 3).Value = 33
$ws.Cells.Item(11, 4).Value = 33
$ws.Cells.Item(11, 8).Value = 21
# Row 12
$ws.Cells.Item(12, 2).Value = 0.6050420168067226
$ws.Cells.Item(12, 3).Value = 72
$ws.Cells.Item(12, 4).Value = 72
$ws.Cells.Item(12, 8).Value = 47
# Row 13
$ws.Cells.Item(13, 2).Value = 0.5217391304347826
$ws.Cells.Item(13, 3).Value = 180
$ws.Cells.Item(13, 4).Value = 180
$ws.Cells.Item(13, 8).Value = 165
# Row 14
$ws.Cells.Item(14, 1).Value = "paint"
$ws.Cells.Item(14, 2).Value = 0.4761904761904762
$ws.Cells.Item(14, 3).Value = 30
$ws.Cells.Item(14, 4).Value = 30
$ws.Cells.Item(14, 8).Value = 33
# Row 15
$ws.Cells.Item(15, 1).Value = "broken"
$ws.Cells.Item(15, 2).Value = 0.4578313253012048
$ws.Cells.Item(15, 3).Value = 38
$ws.Cells.Item(15, 4).Value = 38
$ws.Cells.Item(15, 8).Value = 45
# Row 16
$ws.Cells.Item(16, 1).Value = "apart"
$ws.Cells.Item(16, 2).Value = 0.4421052631578947
$ws.Cells.Item(16, 3).Value = 42
$ws.Cells.Item(16, 4).Value = 42
$ws.Cells.Item(16, 8).Value = 53
# Row 17
$ws.Cells.Item(17, 1).Value = "plastic"
$ws.Cells.Item(17, 2).Value = 0.4330708661417323
$ws.Cells.Item(17, 3).Value = 55
$ws.Cells.Item(17, 4).Value = 55
$ws.Cells.Item(17, 8).Value = 72
# Row 18
$ws.Cells.Item(18, 1).Value = "difficult"
$ws.Cells.Item(18, 2).Value = 0.4044943820224719
$ws.Cells.Item(18, 3).Value = 36
$ws.Cells.Item(18, 4).Value = 36
$ws.Cells.Item(18, 8).Value = 53
# Row 19
$ws.Cells.Item(19, 2).Value = 0.3762376237623762
$ws.Cells.Item(19, 3).Value = 76
$ws.Cells.Item(19, 4).Value = 76
$ws.Cells.Item(19, 8).Value = 126
# Row 20
$ws.Cells.Item(20, 1).Value = "ok"
$ws.Cells.Item(20, 2).Value = 0.375
$ws.Cells.Item(20, 3).Value = 48
$ws.Cells.Item(20, 4).Value = 48
$ws.Cells.Item(20, 8).Value = 80
# Row 21
$ws.Cells.Item(21, 1).Value = "cheap"
$ws.Cells.Item(21, 2).Value = 0.3317535545023697
$ws.Cells.Item(21, 3).Value = 70
$ws.Cells.Item(21, 4).Value = 70
$ws.Cells.Item(21, 8).Value = 141
# Row 22
$ws.Cells.Item(22, 1).Value = "though"
$ws.Cells.Item(22, 2).Value = 0.3076923076923077
$ws.Cells.Item(22, 3).Value = 36
$ws.Cells.Item(22, 4).Value = 36
$ws.Cells.Item(22, 8).Value = 81
# Row 23
$ws.Cells.Item(23, 1).Value = "bit"
$ws.Cells.Item(23, 2).Value = 0.2959183673469388
$ws.Cells.Item(23, 3).Value = 29
$ws.Cells.Item(23, 4).Value = 29
$ws.Cells.Item(23, 8).Value = 69
# Row 24
$ws.Cells.Item(24, 1).Value = "size"
$ws.Cells.Item(24, 2).Value = 0.2680412371134021
$ws.Cells.Item(24, 3).Value = 52
$ws.Cells.Item(24, 4).Value = 52
$ws.Cells.Item(24, 8).Value = 142
# Row 25
$ws.Cells.Item(25, 1).Value = "would"
$ws.Cells.Item(25, 2).Value = 0.2240356083086053
$ws.Cells.Item(25, 3).Value = 151
$ws.Cells.Item(25, 4).Value = 151
$ws.Cells.Item(25, 8).Value = 523
# Row 26
$ws.Cells.Item(26, 1).Value = "money"
$ws.Cells.Item(26, 2).Value = 0.2183544303797468
$ws.Cells.Item(26, 3).Value = 69
$ws.Cells.Item(26, 4).Value = 69
$ws.Cells.Item(26, 8).Value = 247
# Row 27
$ws.Cells.Item(27, 1).Value = "item"
$ws.Cells.Item(27, 2).Value = 0.2173913043478261
$ws.Cells.Item(27, 3).Value = 60
$ws.Cells.Item(27, 4).Value = 60
$ws.Cells.Item(27, 8).Value = 216
# Row 28
$ws.Cells.Item(28, 1).Value = "hard"
$ws.Cells.Item(28, 2).Value = 0.215
$ws.Cells.Item(28, 3).Value = 43
$ws.Cells.Item(28, 4).Value = 43
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 7).Value = $false
$ws.Cells.Item(28, 8).Value = 157
# Row 29
$ws.Cells.Item(29, 1).Value = "work"
$ws.Cells.Item(29, 2).Value = 0.2088607594936709
$ws.Cells.Item(29, 3).Value = 66
$ws.Cells.Item(29, 4).Value = 66
$ws.Cells.Item(29, 8).Value = 250
# Row 30
$ws.Cells.Item(30, 1).Value = "could"
$ws.Cells.Item(30, 2).Value = 0.1974522292993631
$ws.Cells.Item(30, 3).Value = 31
$ws.Cells.Item(30, 4).Value = 31
$ws.Cells.Item(30, 8).Value = 126
# Row 31
$ws.Cells.Item(31, 2).Value = 0.1872246696035242
$ws.Cells.Item(31, 3).Value = 85
$ws.Cells.Item(31, 4).Value = 85
$ws.Cells.Item(31, 8).Value = 369
# Row 32
$ws.Cells.Item(32, 1).Value = "used"
$ws.Cells.Item(32, 2).Value = 0.1657142857142857
$ws.Cells.Item(32, 3).Value = 29
$ws.Cells.Item(32, 4).Value = 29
$ws.Cells.Item(32, 8).Value = 146
# Row 33
$ws.Cells.Item(33, 1).Value = "2"
$ws.Cells.Item(33, 2).Value = 0.1385767790262172
$ws.Cells.Item(33, 3).Value = 37
$ws.Cells.Item(33, 4).Value = 37
$ws.Cells.Item(33, 8).Value = 230
# Row 34
$ws.Cells.Item(34, 1).Value = "better"
$ws.Cells.Item(34, 2).Value = 0.1355140186915888
$ws.Cells.Item(34, 3).Value = 29
$ws.Cells.Item(34, 4).Value = 29
$ws.Cells.Item(34, 8).Value = 185
# Row 35
$ws.Cells.Item(35, 1).Value = "price"
$ws.Cells.Item(35, 2).Value = 0.1350574712643678
$ws.Cells.Item(35, 3).Value = 47
$ws.Cells.Item(35, 4).Value = 47
$ws.Cells.Item(35, 8).Value = 301
# Row 36
$ws.Cells.Item(36, 1).Value = "use"
$ws.Cells.Item(36, 2).Value = 0.1041095890410959
$ws.Cells.Item(36, 3).Value = 38
$ws.Cells.Item(36, 4).Value = 38
$ws.Cells.Item(36, 5).Value = 0
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(36, 7).Value = $false
$ws.Cells.Item(36, 8).Value = 327
# Row 37
$ws.Cells.Item(37, 1).Value = "buy"
$ws.Cells.Item(37, 2).Value = 0.096045197740113
$ws.Cells.Item(37, 3).Value = 34
$ws.Cells.Item(37, 4).Value = 35
$ws.Cells.Item(37, 5).Value = 0.03
$ws.Cells.Item(37, 6).Value = 0.97
$ws.Cells.Item(37, 8).Value = 320

# --- Update existing positive-word rows (3-17): words & stats shifted/changed ---
# Row 3
$ws.Cells.Item(3, 10).Value = "wonderful"
$ws.Cells.Item(3, 11).Value = 0.8392857142857143
$ws.Cells.Item(3, 12).Value = 47
$ws.Cells.Item(3, 13).Value = 47
$ws.Cells.Item(3, 17).Value = 9
# Row 4
$ws.Cells.Item(4, 10).Value = "awesome"
$ws.Cells.Item(4, 11).Value = 0.8153846153846154
$ws.Cells.Item(4, 12).Value = 53
$ws.Cells.Item(4, 13).Value = 53
# Row 5
$ws.Cells.Item(5, 11).Value = 0.6774193548387096
$ws.Cells.Item(5, 12).Value = 63
$ws.Cells.Item(5, 13).Value = 63
$ws.Cells.Item(5, 17).Value = 30
# Row 6
$ws.Cells.Item(6, 11).Value = 0.5849056603773585
$ws.Cells.Item(6, 12).Value = 31
$ws.Cells.Item(6, 13).Value = 31
$ws.Cells.Item(6, 17).Value = 22
# Row 7
$ws.Cells.Item(7, 11).Value = 0.53125
$ws.Cells.Item(7, 12).Value = 34
$ws.Cells.Item(7, 13).Value = 34
$ws.Cells.Item(7, 17).Value = 30
# Row 8
$ws.Cells.Item(8, 11).Value = 0.463768115942029
$ws.Cells.Item(8, 12).Value = 32
$ws.Cells.Item(8, 13).Value = 32
$ws.Cells.Item(8, 17).Value = 37
# Row 9
$ws.Cells.Item(9, 11).Value = 0.3713114754098361
$ws.Cells.Item(9, 12).Value = 453
$ws.Cells.Item(9, 13).Value = 453
$ws.Cells.Item(9, 17).Value = 767
# Row 10
$ws.Cells.Item(10, 11).Value = 0.3247126436781609
$ws.Cells.Item(10, 12).Value = 226
$ws.Cells.Item(10, 16).Value = $true
# Row 11
$ws.Cells.Item(11, 11).Value = 0.3029045643153527
$ws.Cells.Item(11, 12).Value = 146
$ws.Cells.Item(11, 13).Value = 146
$ws.Cells.Item(11, 17).Value = 336
# Row 12
$ws.Cells.Item(12, 10).Value = "loved"
$ws.Cells.Item(12, 11).Value = 0.2171253822629969
$ws.Cells.Item(12, 12).Value = 71
$ws.Cells.Item(12, 13).Value = 71
$ws.Cells.Item(12, 17).Value = 256
# Row 13
$ws.Cells.Item(13, 11).Value = 0.2048192771084337
$ws.Cells.Item(13, 12).Value = 34
$ws.Cells.Item(13, 13).Value = 34
$ws.Cells.Item(13, 17).Value = 132
# Row 14
$ws.Cells.Item(14, 10).Value = "friends"
$ws.Cells.Item(14, 11).Value = 0.201058201058201
$ws.Cells.Item(14, 12).Value = 38
$ws.Cells.Item(14, 13).Value = 38
$ws.Cells.Item(14, 17).Value = 151
# Row 15
$ws.Cells.Item(15, 10).Value = "christmas"
$ws.Cells.Item(15, 11).Value = 0.1285140562248996
$ws.Cells.Item(15, 12).Value = 32
$ws.Cells.Item(15, 13).Value = 32
$ws.Cells.Item(15, 17).Value = 217
# Row 16
$ws.Cells.Item(16, 10).Value = "fun"
$ws.Cells.Item(16, 11).Value = 0.1156879929886065
$ws.Cells.Item(16, 12).Value = 132
$ws.Cells.Item(16, 13).Value = 132
$ws.Cells.Item(16, 17).Value = 1009
# Row 17
$ws.Cells.Item(17, 10).Value = "game"
$ws.Cells.Item(17, 11).Value = 0.04808317089018843
$ws.Cells.Item(17, 12).Value = 74
$ws.Cells.Item(17, 13).Value = 76
$ws.Cells.Item(17, 14).Value = 0.97
$ws.Cells.Item(17, 15).Value = 0.03000000000000003
$ws.Cells.Item(17, 16).Value = $true
$ws.Cells.Item(17, 17).Value = 1465

# --- Row 18 positive side no longer exists: clear J18:Q18 ---
$ws.Range("J18:Q18").ClearContents()

# --- Add new negative-word rows 38-40 ---
# Row 38
$ws.Cells.Item(38, 1).Value = "little"
$ws.Cells.Item(38, 2).Value = 0.07366071428571429
$ws.Cells.Item(38, 3).Value = 33
$ws.Cells.Item(38, 4).Value = 34
$ws.Cells.Item(38, 5).Value = 0.03
$ws.Cells.Item(38, 6).Value = 0.97
$ws.Cells.Item(38, 7).Value = $true
$ws.Cells.Item(38, 8).Value = 415
# Row 39
$ws.Cells.Item(39, 1).Value = "like"
$ws.Cells.Item(39, 2).Value = 0.06930693069306931
$ws.Cells.Item(39, 3).Value = 42
$ws.Cells.Item(39, 4).Value = 44
$ws.Cells.Item(39, 5).Value = 0.05
$ws.Cells.Item(39, 6).Value = 0.95
$ws.Cells.Item(39, 7).Value = $true
$ws.Cells.Item(39, 8).Value = 564
# Row 40
$ws.Cells.Item(40, 1).Value = "one"
$ws.Cells.Item(40, 2).Value = 0.04441624365482234
$ws.Cells.Item(40, 3).Value = 35
$ws.Cells.Item(40, 4).Value = 41
$ws.Cells.Item(40, 5).Value = 0.15
$ws.Cells.Item(40, 6).Value = 0.85
$ws.Cells.Item(40, 7).Value = $true
$ws.Cells.Item(40, 8).Value = 753

# --- Apply the existing bold/centered/bordered "name" style to the new word cells ---
$ws.Cells.Item(2, 1).Copy()
$ws.Range("A38:A40").PasteSpecial(-4122)
$excel.CutCopyMode = $false
